# Automatische test-sync: 2025-07-29 22:11:50
#
# Appends the 20th test-mail (a "Klacht / Probleem" complaint) to the Logs
# sheet, bumps its matching Dashboard category count, and extends the
# dashboard bar chart's category/value ranges to include the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 22
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A22").Value = "Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("B22").Value = "mailmind.test@zohomail.eu"
$logs.Range("C22").Value = "Testmail #20: Ik ben niet tevreden over hoe dit is gegaan."
$logs.Range("D22").Value = "Klacht / Probleem"
$logs.Range("F22").Value = "2025-07-29 22:11:06"
$logs.Range("G22").Value = "Nee"
$logs.Range("H22").Value = "Ja"
$logs.Range("I22").Value = "Nee"
$logs.Range("J22").Value = "Nee"

# Extend the conditional-formatting ranges (each one is currently *2:*21)
# so they keep covering the newly added row 22. Modifying any one rule in
# a FormatConditions collection re-applies the range to the whole group.
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))
$logs.Range("H2:H21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H22"))
$logs.Range("I2:I21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I22"))
$logs.Range("J2:J21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J22"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: append the new category tally row 8
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A8").Value = "Klacht / Probleem"
$dash.Range("B8").Value = 1

# ---------------------------------------------------------------------
# 3. Dashboard chart: extend category/value series ranges from row 7 to 8
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"
